# Append new blood-glucose log readings (2026/02/13 11:08 .. 14:08) as two
# new columns of data starting at row 541, continuing directly after the
# existing last row (540) of Sheet1.
#
# Column A = timestamp ("血糖时间"), Column B = glucose reading ("血糖值").
# Both columns in the existing data are plain text (shared strings), so the
# numeric-looking readings in column B must be forced to text - otherwise
# Excel auto-converts strings like "12.5" into a number. We do that by
# pre-formatting the target B cells as Text ("@") before writing, then
# clearing the (now superfluous) formatting back off again once the values
# are safely stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 541

$newRows = @(
    @("2026/02/13 11:08", "12.5"),
    @("2026/02/13 11:13", "12.4"),
    @("2026/02/13 11:18", "12.3"),
    @("2026/02/13 11:23", "12.2"),
    @("2026/02/13 11:28", "12.1"),
    @("2026/02/13 11:33", "11.6"),
    @("2026/02/13 11:38", "11.8"),
    @("2026/02/13 11:43", "11.6"),
    @("2026/02/13 11:48", "12.1"),
    @("2026/02/13 11:53", "12.0"),
    @("2026/02/13 11:58", "12.0"),
    @("2026/02/13 12:03", "12.4"),
    @("2026/02/13 12:08", "12.4"),
    @("2026/02/13 12:13", "12.0"),
    @("2026/02/13 12:18", "11.8"),
    @("2026/02/13 12:23", "11.9"),
    @("2026/02/13 12:28", "11.7"),
    @("2026/02/13 12:33", "11.6"),
    @("2026/02/13 12:38", "12.3"),
    @("2026/02/13 12:43", "13.7"),
    @("2026/02/13 12:48", "15.6"),
    @("2026/02/13 12:53", "16.9"),
    @("2026/02/13 12:58", "18.0"),
    @("2026/02/13 13:03", "19.3"),
    @("2026/02/13 13:08", "20.3"),
    @("2026/02/13 13:13", "21.8"),
    @("2026/02/13 13:18", "22.7"),
    @("2026/02/13 13:23", "23.4"),
    @("2026/02/13 13:28", "23.8"),
    @("2026/02/13 13:33", "24.8"),
    @("2026/02/13 13:38", "24.5"),
    @("2026/02/13 13:43", "25.0"),
    @("2026/02/13 13:48", "25.0"),
    @("2026/02/13 13:53", "24.8"),
    @("2026/02/13 13:58", "24.7"),
    @("2026/02/13 14:03", "25.0"),
    @("2026/02/13 14:08", "24.8")
)

$endRow = $startRow + $newRows.Count - 1

# Pre-mark the value column as Text so the numeric-looking readings are
# stored as literal strings (matching the existing column B's shared-string
# data) instead of being auto-coerced to numbers.
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $newRows[$i][0]
    $ws.Range("B$r").Value = $newRows[$i][1]
}

# Drop the temporary Text formatting again now that the values are locked in
# as strings, restoring the cells to the sheet's default (unstyled) look.
$ws.Range("A$startRow`:B$endRow").ClearFormats()
